# Prize pool values for the "20% paid" leaderboard payouts table (B2:AC37)
# were reworked from fractional percentages into flat cent/point amounts.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("20% paid")

$rowValues = @(1550, 1600, 1650, 1700, 1750, 1800, 1850, 1900, 1950, 2000, 2050, 2100, 2150, 2200, 2250, 2300, 2400, 2400, 2600, 2700, 2800, 2900, 3100, 3300, 3500, 4000, 7000, 10000)
for ($i = 0; $i -lt $rowValues.Count; $i++) { $ws.Cells.Item(2, $i + 2).Value = $rowValues[$i] }
$rowValues = @(1075, 1100, 1125, 1150, 1175, 1200, 1225, 1250, 1275, 1300, 1325, 1350, 1375, 1400, 1425, 1450, 1500, 1550, 1600, 1650, 1700, 1750, 1900, 2000, 2200, 2500, 3000)
for ($i = 0; $i -lt $rowValues.Count; $i++) { $ws.Cells.Item(3, $i + 2).Value = $rowValues[$i] }
$rowValues = @(640, 640, 660, 660, 670, 680, 690, 690, 700, 720, 720, 740, 740, 750, 760, 770, 780, 790, 810, 820, 840, 890, 980, 1200, 1500, 2000)
for ($i = 0; $i -lt $rowValues.Count; $i++) { $ws.Cells.Item(4, $i + 2).Value = $rowValues[$i] }
$rowValues = @(540, 540, 560, 560, 570, 580, 590, 590, 600, 620, 620, 640, 640, 650, 660, 670, 680, 690, 710, 720, 750, 790, 880, 900, 1100, 1500)
for ($i = 0; $i -lt $rowValues.Count; $i++) { $ws.Cells.Item(5, $i + 2).Value = $rowValues[$i] }
$rowValues = @(440, 440, 460, 460, 470, 480, 490, 490, 500, 520, 520, 540, 540, 550, 560, 570, 580, 590, 610, 620, 650, 690, 780, 800, 900)
for ($i = 0; $i -lt $rowValues.Count; $i++) { $ws.Cells.Item(6, $i + 2).Value = $rowValues[$i] }
$rowValues = @(340, 340, 350, 350, 360, 370, 380, 380, 400, 420, 420, 440, 440, 450, 460, 470, 480, 490, 510, 520, 550, 590, 680, 700, 800)
for ($i = 0; $i -lt $rowValues.Count; $i++) { $ws.Cells.Item(7, $i + 2).Value = $rowValues[$i] }
$rowValues = @(230, 230, 240, 240, 250, 260, 270, 270, 280, 300, 300, 320, 320, 330, 340, 350, 360, 370, 390, 400, 430, 470, 550, 600)
for ($i = 0; $i -lt $rowValues.Count; $i++) { $ws.Cells.Item(8, $i + 2).Value = $rowValues[$i] }
$rowValues = @(130, 130, 140, 140, 140, 150, 160, 160, 170, 180, 180, 190, 190, 200, 210, 220, 230, 240, 260, 270, 290, 350, 440, 500)
for ($i = 0; $i -lt $rowValues.Count; $i++) { $ws.Cells.Item(9, $i + 2).Value = $rowValues[$i] }
$rowValues = @(95, 95, 100, 100, 100, 110, 120, 120, 130, 150, 150, 160, 160, 170, 180, 190, 200, 210, 230, 240, 260, 290, 370)
for ($i = 0; $i -lt $rowValues.Count; $i++) { $ws.Cells.Item(10, $i + 2).Value = $rowValues[$i] }
$rowValues = @(75, 75, 80, 80, 80, 90, 95, 95, 105, 110, 110, 115, 115, 120, 125, 130, 140, 145, 150, 160, 180, 230, 320)
for ($i = 0; $i -lt $rowValues.Count; $i++) { $ws.Cells.Item(11, $i + 2).Value = $rowValues[$i] }
$rowValues = @(75, 75, 80, 80, 80, 90, 95, 95, 105, 110, 110, 115, 115, 120, 125, 130, 140, 145, 155, 160, 180, 210)
for ($i = 0; $i -lt $rowValues.Count; $i++) { $ws.Cells.Item(12, $i + 2).Value = $rowValues[$i] }
$rowValues = @(35, 35, 35, 45, 45, 50, 50, 55, 60, 60, 60, 65, 70, 75, 75, 80, 80, 95, 105, 115, 130)
for ($i = 0; $i -lt $rowValues.Count; $i++) { $ws.Cells.Item(13, $i + 2).Value = $rowValues[$i] }
$rowValues = @(30, 30, 30, 35, 35, 40, 40, 45, 50, 50, 50, 50, 55, 60, 65, 70, 70, 75, 90, 105)
for ($i = 0; $i -lt $rowValues.Count; $i++) { $ws.Cells.Item(14, $i + 2).Value = $rowValues[$i] }
$rowValues = @(20, 20, 20, 25, 25, 30, 30, 35, 40, 40, 40, 40, 45, 50, 55, 55, 55, 60, 75)
for ($i = 0; $i -lt $rowValues.Count; $i++) { $ws.Cells.Item(15, $i + 2).Value = $rowValues[$i] }
$rowValues = @(19, 19, 19, 20, 20, 20, 20, 25, 30, 30, 30, 30, 35, 40, 45, 50, 50, 55)
for ($i = 0; $i -lt $rowValues.Count; $i++) { $ws.Cells.Item(16, $i + 2).Value = $rowValues[$i] }
$rowValues = @(16, 16, 16, 18, 18, 18, 18, 21, 25, 25, 25, 25, 30, 35, 40, 45, 45, 55)
for ($i = 0; $i -lt $rowValues.Count; $i++) { $ws.Cells.Item(17, $i + 2).Value = $rowValues[$i] }
$rowValues = @(15, 15, 15, 16, 16, 17, 17, 18, 21, 21, 21, 21, 25, 31, 35, 40, 45)
for ($i = 0; $i -lt $rowValues.Count; $i++) { $ws.Cells.Item(18, $i + 2).Value = $rowValues[$i] }
$rowValues = @(14, 14, 14, 14, 14, 15, 15, 17, 18, 18, 19, 19, 22, 26, 31, 33)
for ($i = 0; $i -lt $rowValues.Count; $i++) { $ws.Cells.Item(19, $i + 2).Value = $rowValues[$i] }
$rowValues = @(13, 13, 13, 13, 13, 13, 13, 15, 15, 16, 17, 17, 19, 19, 23)
for ($i = 0; $i -lt $rowValues.Count; $i++) { $ws.Cells.Item(20, $i + 2).Value = $rowValues[$i] }
$rowValues = @(12, 12, 12, 12, 12, 12, 12, 12, 12, 13, 15, 16, 17, 17)
for ($i = 0; $i -lt $rowValues.Count; $i++) { $ws.Cells.Item(21, $i + 2).Value = $rowValues[$i] }
$rowValues = @(11, 11, 11, 11, 11, 11, 11, 11, 11, 12, 14, 15, 16)
for ($i = 0; $i -lt $rowValues.Count; $i++) { $ws.Cells.Item(22, $i + 2).Value = $rowValues[$i] }
$rowValues = @(10, 10, 10, 10, 10, 11, 10, 10, 10, 12, 13, 14)
for ($i = 0; $i -lt $rowValues.Count; $i++) { $ws.Cells.Item(23, $i + 2).Value = $rowValues[$i] }
$rowValues = @(9, 9, 9, 9, 9, 9, 9, 9, 9, 11, 13)
for ($i = 0; $i -lt $rowValues.Count; $i++) { $ws.Cells.Item(24, $i + 2).Value = $rowValues[$i] }
$rowValues = @(8, 8, 8, 8, 8, 8, 8, 9, 9, 11)
for ($i = 0; $i -lt $rowValues.Count; $i++) { $ws.Cells.Item(25, $i + 2).Value = $rowValues[$i] }
$rowValues = @(7, 7, 7, 7, 7, 7, 8, 9, 9)
for ($i = 0; $i -lt $rowValues.Count; $i++) { $ws.Cells.Item(26, $i + 2).Value = $rowValues[$i] }
$rowValues = @(6, 7, 7, 7, 7, 7, 8, 8, 9)
for ($i = 0; $i -lt $rowValues.Count; $i++) { $ws.Cells.Item(27, $i + 2).Value = $rowValues[$i] }
$rowValues = @(5, 6, 6, 6, 7, 7, 7, 8)
for ($i = 0; $i -lt $rowValues.Count; $i++) { $ws.Cells.Item(28, $i + 2).Value = $rowValues[$i] }
$rowValues = @(5, 6, 6, 6, 6, 6, 7, 8)
for ($i = 0; $i -lt $rowValues.Count; $i++) { $ws.Cells.Item(29, $i + 2).Value = $rowValues[$i] }
$rowValues = @(5, 6, 6, 6, 6, 6, 7)
for ($i = 0; $i -lt $rowValues.Count; $i++) { $ws.Cells.Item(30, $i + 2).Value = $rowValues[$i] }
$rowValues = @(5, 5, 5, 5, 6, 6)
for ($i = 0; $i -lt $rowValues.Count; $i++) { $ws.Cells.Item(31, $i + 2).Value = $rowValues[$i] }
$rowValues = @(5, 5, 5, 5, 6)
for ($i = 0; $i -lt $rowValues.Count; $i++) { $ws.Cells.Item(32, $i + 2).Value = $rowValues[$i] }
$rowValues = @(5, 5, 5, 5)
for ($i = 0; $i -lt $rowValues.Count; $i++) { $ws.Cells.Item(33, $i + 2).Value = $rowValues[$i] }
$rowValues = @(5, 4, 4)
for ($i = 0; $i -lt $rowValues.Count; $i++) { $ws.Cells.Item(34, $i + 2).Value = $rowValues[$i] }
$rowValues = @(5, 4)
for ($i = 0; $i -lt $rowValues.Count; $i++) { $ws.Cells.Item(35, $i + 2).Value = $rowValues[$i] }
$rowValues = @(5)
for ($i = 0; $i -lt $rowValues.Count; $i++) { $ws.Cells.Item(36, $i + 2).Value = $rowValues[$i] }
$rowValues = @(5)
for ($i = 0; $i -lt $rowValues.Count; $i++) { $ws.Cells.Item(37, $i + 2).Value = $rowValues[$i] }
